$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(20, 8).Value = 5910.4
$ws.Cells.Item(20, 9).Value = 1388
$ws.Cells.Item(20, 10).Value = 24000
$ws.Cells.Item(20, 11).Value = 1388
$ws.Cells.Item(20, 12).Value = 24000
$ws.Cells.Item(20, 13).Value = -1158
$ws.Cells.Item(20, 14).Value = -24460

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(35, 8).Value = 5910.4
$ws.Cells.Item(35, 9).Value = 1388
$ws.Cells.Item(35, 10).Value = 24000
$ws.Cells.Item(35, 11).Value = 1388
$ws.Cells.Item(35, 12).Value = 24000
$ws.Cells.Item(35, 13).Value = -1009
$ws.Cells.Item(35, 14).Value = -24758

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 2115.4
$ws.Cells.Item(43, 9).Value = 5000
$ws.Cells.Item(43, 10).Value = 1394.25
$ws.Cells.Item(43, 11).Value = 5000
$ws.Cells.Item(43, 12).Value = 1394.25
$ws.Cells.Item(43, 13).Value = -4931
$ws.Cells.Item(43, 14).Value = -1532.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 1951
$ws.Cells.Item(141, 9).Value = 1000
$ws.Cells.Item(141, 10).Value = 2902
$ws.Cells.Item(141, 11).Value = 3000
$ws.Cells.Item(141, 12).Value = 8706
$ws.Cells.Item(141, 13).Value = 2180
$ws.Cells.Item(141, 14).Value = -19066

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(19, 8).Value = 5602.6665
$ws.Cells.Item(19, 9).Value = 5399
$ws.Cells.Item(19, 10).Value = 6010
$ws.Cells.Item(19, 11).Value = 5399
$ws.Cells.Item(19, 12).Value = 6010
$ws.Cells.Item(19, 13).Value = -5170
$ws.Cells.Item(19, 14).Value = -6468

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(38, 8).Value = 200
$ws.Cells.Item(38, 9).Value = 200
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 200
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 13).Value = 267

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(19, 8).Value = 9000
$ws.Cells.Item(19, 9).Value = 9000
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 9000
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = -8827

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 6548.143
$ws.Cells.Item(86, 9).Value = 11581.667
$ws.Cells.Item(86, 10).Value = 2773
$ws.Cells.Item(86, 11).Value = 11581.667
$ws.Cells.Item(86, 12).Value = 2773
$ws.Cells.Item(86, 13).Value = -10458.667
$ws.Cells.Item(86, 14).Value = -5019

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 6548.143
$ws.Cells.Item(89, 9).Value = 11581.667
$ws.Cells.Item(89, 10).Value = 2773
$ws.Cells.Item(89, 11).Value = 57908.335
$ws.Cells.Item(89, 12).Value = 13865
$ws.Cells.Item(89, 13).Value = -52292.335
$ws.Cells.Item(89, 14).Value = -25097

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).ClearContents()
$ws.Cells.Item(15, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(42, 8).Value = 6287.3335
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 6287.3335
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 6287.3335
$ws.Cells.Item(42, 13).ClearContents()
$ws.Cells.Item(42, 14).Value = -7473.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 3476364.5
$ws.Cells.Item(62, 9).Value = 6176450.5
$ws.Cells.Item(62, 10).Value = 4825.2856
$ws.Cells.Item(62, 11).Value = 6176450.5
$ws.Cells.Item(62, 12).Value = 4825.2856
$ws.Cells.Item(62, 13).Value = -6175826.5
$ws.Cells.Item(62, 14).Value = -6073.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 3476364.5
$ws.Cells.Item(65, 9).Value = 6176450.5
$ws.Cells.Item(65, 10).Value = 4825.2856
$ws.Cells.Item(65, 11).Value = 30882252.5
$ws.Cells.Item(65, 12).Value = 24126.428
$ws.Cells.Item(65, 13).Value = -30879132.5
$ws.Cells.Item(65, 14).Value = -30366.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 31251990
$ws.Cells.Item(134, 9).Value = 55556572
$ws.Cells.Item(134, 10).Value = 3244.7144
$ws.Cells.Item(134, 11).Value = 166669716
$ws.Cells.Item(134, 12).Value = 9734.143199999999
$ws.Cells.Item(134, 13).Value = -166667181
$ws.Cells.Item(134, 14).Value = -14804.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 3851.2856
$ws.Cells.Item(3, 9).Value = 1992
$ws.Cells.Item(3, 10).Value = 7198
$ws.Cells.Item(3, 11).Value = 5976
$ws.Cells.Item(3, 12).Value = 21594
$ws.Cells.Item(3, 13).Value = -5864
$ws.Cells.Item(3, 14).Value = -21818

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(36, 8).Value = 672.38464
$ws.Cells.Item(36, 9).Value = 208.4
$ws.Cells.Item(36, 10).Value = 962.375
$ws.Cells.Item(36, 11).Value = 625.2
$ws.Cells.Item(36, 12).Value = 2887.125
$ws.Cells.Item(36, 13).Value = -456.2
$ws.Cells.Item(36, 14).Value = -3225.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(41, 8).Value = 398.75
$ws.Cells.Item(41, 9).Value = 275
$ws.Cells.Item(41, 10).Value = 440
$ws.Cells.Item(41, 11).Value = 825
$ws.Cells.Item(41, 12).Value = 1320
$ws.Cells.Item(41, 13).Value = -487
$ws.Cells.Item(41, 14).Value = -1996

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(42, 8).Value = 5000
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 5000
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 15000
$ws.Cells.Item(42, 14).Value = -16068

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(81, 8).Value = 2417.0833
$ws.Cells.Item(81, 9).Value = 1496.6666
$ws.Cells.Item(81, 10).Value = 2723.889
$ws.Cells.Item(81, 11).Value = 4489.9998
$ws.Cells.Item(81, 12).Value = 8171.667
$ws.Cells.Item(81, 13).Value = -3366.9998
$ws.Cells.Item(81, 14).Value = -10417.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(84, 8).Value = 2417.0833
$ws.Cells.Item(84, 9).Value = 1496.6666
$ws.Cells.Item(84, 10).Value = 2723.889
$ws.Cells.Item(84, 11).Value = 13469.9994
$ws.Cells.Item(84, 12).Value = 24515.001
$ws.Cells.Item(84, 13).Value = -7853.999400000001
$ws.Cells.Item(84, 14).Value = -35747.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5140.6665
$ws.Cells.Item(70, 9).Value = 5680.269
$ws.Cells.Item(70, 10).Value = 4061.4614
$ws.Cells.Item(70, 11).Value = 5680.269
$ws.Cells.Item(70, 12).Value = 4061.4614
$ws.Cells.Item(70, 13).Value = -5410.269
$ws.Cells.Item(70, 14).Value = -4601.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 5140.6665
$ws.Cells.Item(73, 9).Value = 5680.269
$ws.Cells.Item(73, 10).Value = 4061.4614
$ws.Cells.Item(73, 11).Value = 5680.269
$ws.Cells.Item(73, 12).Value = 4061.4614
$ws.Cells.Item(73, 13).Value = -4744.269
$ws.Cells.Item(73, 14).Value = -5933.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(93, 8).Value = 31800
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 31800
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 31800
$ws.Cells.Item(93, 14).Value = -35544

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 779.4262
$ws.Cells.Item(46, 9).Value = 698.5526
$ws.Cells.Item(46, 10).Value = 913.04346
$ws.Cells.Item(46, 11).Value = 698.5526
$ws.Cells.Item(46, 12).Value = 913.04346
$ws.Cells.Item(46, 13).Value = -510.5526
$ws.Cells.Item(46, 14).Value = -1289.04346

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(47, 8).Value = 9955
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 9955
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 12).Value = 9955
$ws.Cells.Item(47, 14).Value = -10935

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(52, 8).Value = 9955
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 10).Value = 9955
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 12).Value = 9955
$ws.Cells.Item(52, 14).Value = -10421

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(53, 8).Value = 12566.667
$ws.Cells.Item(53, 9).Value = 500
$ws.Cells.Item(53, 10).Value = 18600
$ws.Cells.Item(53, 11).Value = 500
$ws.Cells.Item(53, 12).Value = 18600
$ws.Cells.Item(53, 13).Value = 18
$ws.Cells.Item(53, 14).Value = -19636

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(95, 8).Value = 13789.8
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = 13789.8
$ws.Cells.Item(95, 11).Value = 0
$ws.Cells.Item(95, 12).Value = 13789.8
$ws.Cells.Item(95, 14).Value = -19281.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(42, 8).Value = 3000
$ws.Cells.Item(42, 9).Value = 3000
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 3000
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = -2622
$ws.Cells.Item(42, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(97, 8).Value = 21000
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 21000
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 21000
$ws.Cells.Item(97, 14).Value = -22982
